# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" tab (fund holdings detail) right after the "总计"
# (summary) tab, and updates the "总计" tab with a new leading row for the
# 2022-Q4 quarter (every other quarter's summary row simply slides down by
# one position, values unchanged).
# ---------------------------------------------------------------------------

function Set-SmartValue($range, $val) {
    # Excel infers type from the literal the same way typing into the grid
    # would: a bare numeric-looking string becomes a *number* cell (and
    # loses leading zeros, e.g. fund code "010365"). Prefixing with a
    # leading apostrophe is the standard Excel trick to force text storage
    # (quoted-number / "number stored as text"), which is what we need for
    # fund codes and the percentage-ish text columns that are stored as
    # text in the source data.
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $val
    } else {
        $range.Value = $val
    }
}

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Sheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Duplicate the "2022-Q3" detail sheet (same header / column styling as
#    every other quarterly detail tab) right after "总计", rename it to
#    "2022-Q4", then overwrite its data rows with the 2022-Q4 figures.
# ---------------------------------------------------------------------------
$sourceSheet = $wb.Sheets.Item("2022-Q3")
$sourceSheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Sheets.Item(2)
$q4Sheet.Name = "2022-Q4"

$q4Data = @(
  @('008515', '富兰克林国海基本面优选混合', '12.93', '88.93', '7.21', '0.9323', 3),
  @('011152', '富兰克林国海兴海回报混合', '16.09', '87.64', '4.61', '0.7417', 8),
  @('010365', '鹏华港股通中证香港银行投资指数（LOF）C', '6.46', '94.47', '2.21', '0.1428', 10),
  @('513690', '博时恒生港股通高股息率ETF', '5.20', '98.04', '2.67', '0.1388', 5),
  @('501025', '鹏华港股通中证香港银行投资指数（LOF）A', '2.38', '94.47', '2.21', '0.0526', 10),
  @('501305', '汇添富中证港股通高股息投资指数（LOF）A', '0.89', '91.24', '3.26', '0.0290', 9),
  @('159726', '华夏恒生中国内地企业高股息率ETF', '0.84', '98.34', '2.93', '0.0246', 4),
  @('007751', '景顺长城中证沪港深红利成长低波动指数A', '0.69', '91.50', '3.33', '0.0230', 4),
  @('513530', '华泰柏瑞中证港股通高股息投资ETF（QDII）', '0.63', '96.34', '3.43', '0.0216', 9),
  @('006810', '泰康港股通中证香港银行投资指数C', '0.58', '94.66', '2.20', '0.0128', 10),
  @('006809', '泰康港股通中证香港银行投资指数A', '0.55', '94.66', '2.20', '0.0121', 10),
  @('501306', '汇添富中证港股通高股息投资指数（LOF）C', '0.23', '91.24', '3.26', '0.0075', 9),
  @('005702', '恒生前海港股通高股息低波动指数', '0.23', '94.47', '2.52', '0.0058', 7),
  @('007760', '景顺长城中证沪港深红利成长低波动指数C', '0.07', '91.50', '3.33', '0.0023', 4)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = 2 + $i
    $row = $q4Data[$i]
    Set-SmartValue $q4Sheet.Cells.Item($r, 2) $row[0]
    Set-SmartValue $q4Sheet.Cells.Item($r, 3) $row[1]
    Set-SmartValue $q4Sheet.Cells.Item($r, 4) $row[2]
    Set-SmartValue $q4Sheet.Cells.Item($r, 5) $row[3]
    Set-SmartValue $q4Sheet.Cells.Item($r, 6) $row[4]
    Set-SmartValue $q4Sheet.Cells.Item($r, 7) $row[5]
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2) "总计" tab: insert a new row 2 (carrying formatting down from the
#    existing rows) and then rewrite every data row explicitly with the
#    final values, so the quarter labels / counts / market values all line
#    up correctly regardless of how the row-insert shifted formats.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalData = @(
  @(0, "2022-Q4", 14, 2.15),
  @(1, "2022-Q3", 14, 1.74),
  @(2, "2022-Q2", 7, 1.03),
  @(3, "2022-Q1", 16, 1.37),
  @(4, "2021-Q4", 5, 0.91),
  @(5, "2021-Q3", 7, 1.01),
  @(6, "2021-Q2", 4, 1.03),
  @(7, "2021-Q1", 7, 1.18),
  @(8, "2020-Q4", 1, 0.7)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = 2 + $i
    $row = $totalData[$i]
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}
